$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 2-6 (2005年-2009年), shifting 2010年-2013年 (rows 7-10) up to rows 2-5
$ws.Range("A2:I6").Delete()
